$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last student (ID 211892) originally had 5 duplicate rows (315-319), one
# per Group (B2A, B2B, B2C, B2D, B2E). Keep only a single row for the
# student, now pointing at Group B2B, clear out the now-unused duplicate
# rows' contents, and remove the trailing blank row.

# 1) Row 315 now reflects Group B2B instead of Group B2A.
$ws.Range("D315").Value = "B2B"

# 2) The old duplicate rows (formerly B2C/B2D/B2E) are cleared out, leaving
#    blank (but still present) rows.
$ws.Range("A316:E318").ClearContents()

# 3) The last duplicate row is removed entirely, shrinking the used range.
$ws.Rows(319).Delete()

# Restore the cursor/selection to the edited row.
$ws.Rows("315:315").Select() | Out-Null
